$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C27").Value = "[name=`"Zofia`"]  No! And don’t call me 'Auntie!' `n"
$ws.Range("C33").Value = "[name=`"Zofia`"]  Maria... Call me 'Auntie'... one more time... I’ll... mess ya upp...`n"
$ws.Range("C38").Value = "[name=`"Zofia`"]  Who’s yer 'Auntie?!'`n"
$ws.Range("C53").Value = "[name=`"Old Craftsman`"]  But that Leithanien lass. She actually deserves the title 'Black'. She’s a monster, no mistake. `n"
$ws.Range("C61").Value = "[name=`"Old Knight`"]  Her opponent was the 'Fissure'. That reminds me, did he enter the Major this year?`n"
$ws.Range("C82").Value = "[name=`"Old Knight`"]  Shut up. That’s exactly why 'Fissure' defeated Zofia.`n"
$ws.Range("C127").Value = "[name=`"Corporate Employee`"]  I see, so 'Whislash' is Ms. Maria’s coach. All the better.`n"
$ws.Range("C176").Value = "[name=`"??? `"]  ...Ingra? 'Brassrust' Ingra? The National Council let him go again?  `n"
$ws.Range("C188").Value = "[name=`"Greatmouth Mob`"]  That’s right. He’s the one who been taken to court countless times for crippling anyone who dares cross him, and thanks to the coffers of his countless sponsors has been ruled innocent over and over by the National Council! 'Brassrust' Ingra!`n"
$ws.Range("C197").Value = "[name=`"Greatmouth Mob`"]  Will all this change because of this one match? Will 'Brassrust' bash in this pretty little face of hers? `n"
